# Add two lines of body text to the "Объект 2" placeholder (idx=1) on
# slide 7 ("Проблемы и сложности" / Problems and difficulties), which
# currently is empty:
#   - "Большое кол-во объектов"
#   - "Особенности Pygame"
# Each line becomes its own paragraph, matching the target diff.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(7)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# First paragraph.
$tr.InsertBefore("Большое кол-во объектов")
$tr.LanguageID = "ru-RU"

# Second paragraph (separated with a hard paragraph break "`r").
$tr.InsertAfter("`rОсобенности Pygame")
$tr.LanguageID = "ru-RU"
